# Generate Report for Handoff
# The localization status moved from "In Translation" to "Ready for handoff",
# and the handoff timestamps were refreshed accordingly across the
# Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-11-09 06:38:48"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-11-09 06:38:35"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-11-09 06:38:48"

# The wider "Ready for handoff" text no longer fits the old column width,
# so widen the affected status columns to keep the report readable.
$overview.Columns.Item(5).ColumnWidth = 17.2159881591797
$overview.Columns.Item(6).ColumnWidth = 17.2159881591797
$zhcn.Columns.Item(3).ColumnWidth = 17.2159881591797
$dede.Columns.Item(3).ColumnWidth = 17.2159881591797
